$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("D2").ClearContents()
$ws.Range("C2").Value = 41.226569807504134

$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

$ws.Range("B1:E3").Select()
